$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.665.44"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "2.207.11"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'253.93"
$ws.Range("E5").Value = "  +3.18%  "

$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("D7").Value = "'74.62"
$ws.Range("E7").Value = "  -2.24%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  -4.60%  "

$ws.Range("D10").Value = "'40.82"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("D11").Value = "'0.0915"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "'6.87"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").Value = "2.541.19"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").Value = "'14.23"
$ws.Range("E15").Value = "  -3.02%  "

$ws.Range("D16").Value = "2.196.03"
$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("D17").Value = "'0.774"
$ws.Range("E17").Value = "  -4.88%  "

$ws.Range("D18").Value = "42.585.62"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").Value = "'71.00"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").Value = "'5.94"
$ws.Range("E21").Value = "  -0.92%  "

$ws.Range("D22").Value = "'228.34"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("D24").Value = "'9.39"
$ws.Range("E24").Value = "  -8.50%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "'10.54"
$ws.Range("E26").Value = "  -3.29%  "

$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").Value = "'39.53"
$ws.Range("E27").Value = "  +3.61%  "

$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").Value = "'3.37"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("E29").Value = "  +2.89%  "

$ws.Range("E30").Value = "  -3.14%  "

$ws.Range("D31").Value = "'172.72"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").Value = "'20.16"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("D33").Value = "'0.0832"
$ws.Range("E33").Value = "  +4.50%  "

$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").Value = "'0.109"
$ws.Range("E36").Value = "  -3.93%  "

$ws.Range("D37").Value = "'0.0342"
$ws.Range("E37").Value = "  +4.52%  "

$ws.Range("D38").Value = "'4.27"
$ws.Range("E38").Value = "  -1.96%  "

$ws.Range("D39").Value = "'12.29"
$ws.Range("E39").Value = "  -5.34%  "

$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = "  +17.81%  "

$ws.Range("D42").Value = "'5.25"
$ws.Range("E42").Value = "  -5.75%  "

$ws.Range("D43").Value = "'59.75"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("E44").Value = "  -3.28%  "

$ws.Range("D45").Value = "'102.43"
$ws.Range("E45").Value = "  -2.88%  "

$ws.Range("D46").Value = "'8.34"
$ws.Range("E46").Value = "  -3.68%  "

$ws.Range("D47").Value = "'0.0977"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").Value = "'0.457"
$ws.Range("E48").Value = "  +3.85%  "

$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("E50").Value = "  -1.40%  "

$ws.Range("D51").Value = "2.436.15"
$ws.Range("E51").Value = "  -0.16%  "
